$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B57").Value = "LSBF"
$ws.Range("A58").Value = "Vlaamse Gevarentaal"
$ws.Range("B58").Value = "VGT"
$ws.Range("A59").Value = "Système des signes international"
$ws.Range("B59").Value = "SSI"
